# Generate Report for Handoff
# - The previously "Ready for handoff" dependency pair (3cf4ae3f....png /
#   dd506806....png) is gone from this run.
# - A brand-new source markdown file (0ea69e5b-be3c-480b-818e-fd0f27e57300.md)
#   replaces 158279b2-0d2c-47c2-98ad-ab23967059fc.md, and a second new
#   markdown file (b363c423-27c9-4873-8c27-90ce290b92c5.md) is now tracked
#   too, each with its own freshly generated handoff .xlf package.
# - .localization-config keeps trailing as the "Not to be localized" row.

$wb = $excel.ActiveWorkbook

$e2eBase    = "https://github.com/OpenLocalizationTest/oltest/blob/60167520f76feec5d37c728c9896190337ce0d7c/e2e/"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/60167520f76feec5d37c728c9896190337ce0d7c/.localization-config"
$zhcnXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c624dac73d7fc5657f41f5a231376cdd3d0fa65e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$dedeXlfBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efa9e5a87a1e12f3599bc8affcfc8b2c528edf91/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$md1 = "0ea69e5b-be3c-480b-818e-fd0f27e57300.md"
$md2 = "b363c423-27c9-4873-8c27-90ce290b92c5.md"
$cfg = ".localization-config"

$xlf1zh = "0ea69e5b-be3c-480b-818e-fd0f27e57300.a8f0fff59071464705fcd1b95a3c9b30a811cc4f.zh-cn.xlf"
$xlf2zh = "b363c423-27c9-4873-8c27-90ce290b92c5.8c5ddec2e036b086f2dffe31798fa2ec1d528d5b.zh-cn.xlf"
$xlf1de = "0ea69e5b-be3c-480b-818e-fd0f27e57300.a8f0fff59071464705fcd1b95a3c9b30a811cc4f.de-de.xlf"
$xlf2de = "b363c423-27c9-4873-8c27-90ce290b92c5.8c5ddec2e036b086f2dffe31798fa2ec1d528d5b.de-de.xlf"

$readyForHandoff   = "Ready for handoff"
$notToBeLocalized  = "Not to be localized"
$include           = "Include"
$ignored           = "Ignored"
$zeroDate          = "0001-01-01 00:00:00"
$zhTimestamp       = "2016-03-09 08:40:55"
$deTimestamp       = "2016-03-09 08:41:00"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = $readyForHandoff
$ws.Range("C2").Value = $readyForHandoff

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = $readyForHandoff
$ws.Range("C3").Value = $readyForHandoff

$ws.Range("A4").Value = $cfg
$ws.Range("B4").Value = $notToBeLocalized
$ws.Range("C4").Value = $notToBeLocalized

$ws.Rows("5").Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + $md1), "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + $md2), "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $cfg)

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = $readyForHandoff
$ws.Range("C2").Value = $xlf1zh
$ws.Range("D2").Value = $zhTimestamp
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = $zeroDate
$ws.Range("H2").Value = $include
$ws.Range("I2").Value = ""

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = $readyForHandoff
$ws.Range("C3").Value = $xlf2zh
$ws.Range("D3").Value = $zhTimestamp
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = $zeroDate
$ws.Range("H3").Value = $include
$ws.Range("I3").Value = ""

$ws.Range("A4").Value = $cfg
$ws.Range("B4").Value = $notToBeLocalized
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = $zeroDate
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = $zeroDate
$ws.Range("H4").Value = $ignored
$ws.Range("I4").Value = ""

$ws.Rows("5").Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + $md1), "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("C2"), ($zhcnXlfBase + $xlf1zh), "", "", $xlf1zh)
$ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + $md2), "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("C3"), ($zhcnXlfBase + $xlf2zh), "", "", $xlf2zh)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $cfg)

# ---------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $md1
$ws.Range("B2").Value = $readyForHandoff
$ws.Range("C2").Value = $xlf1de
$ws.Range("D2").Value = $deTimestamp
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = $zeroDate
$ws.Range("H2").Value = $include
$ws.Range("I2").Value = ""

$ws.Range("A3").Value = $md2
$ws.Range("B3").Value = $readyForHandoff
$ws.Range("C3").Value = $xlf2de
$ws.Range("D3").Value = $deTimestamp
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = $zeroDate
$ws.Range("H3").Value = $include
$ws.Range("I3").Value = ""

$ws.Range("A4").Value = $cfg
$ws.Range("B4").Value = $notToBeLocalized
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = $zeroDate
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = $zeroDate
$ws.Range("H4").Value = $ignored
$ws.Range("I4").Value = ""

$ws.Rows("5").Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + $md1), "", "", $md1)
$ws.Hyperlinks.Add($ws.Range("C2"), ($dedeXlfBase + $xlf1de), "", "", $xlf1de)
$ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + $md2), "", "", $md2)
$ws.Hyperlinks.Add($ws.Range("C3"), ($dedeXlfBase + $xlf2de), "", "", $xlf2de)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $cfg)

$wb.Worksheets.Item("Overview").Activate()
